# "Updated the Column EUC in EUC_UIP.docx"
#
# The "EUC" column (column H) on the "Inc Matrix" sheet had several rows
# still marked "Not Started" by mistake. Correct them:
#   - UC 8  (row 9)  and UC 21 (row 22) -> "Double Check" (new status value)
#   - UC 13-15 (rows 14-16) and UC 19-20 (rows 20-21) -> "Done"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H9").Value  = "Double Check"
$ws.Range("H14").Value = "Done"
$ws.Range("H15").Value = "Done"
$ws.Range("H16").Value = "Done"
$ws.Range("H20").Value = "Done"
$ws.Range("H21").Value = "Done"
$ws.Range("H22").Value = "Double Check"

# Column H widened slightly to fit the new, longer "Double Check" text.
$ws.Columns.Item(8).ColumnWidth = 12.5

# Cursor/selection ended up on J5 when the edit was saved.
$ws.Range("J5").Select()
